$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 previously had only A13:C13; add a new D13 cell with a fresh time
# estimate. Set this BEFORE the row 9 update so the new shared strings are
# appended to sharedStrings.xml in the same order as in the target workbook
# ("~5:00" before "~10:00 (...)").
$ws.Range("D13").Value = "~5:00"

# Row 9's time estimate changes from "9:38 (...)" to "~10:00 (...)".
$ws.Range("D9").Value = "~10:00 (using most naïve algorithm; can definitely be improved by only looking at stations from neighboring states)"

# The wrapped text in D9 is now one line taller, so the row grows from
# 71.25 to 85.5 points.
$ws.Rows.Item(9).RowHeight = 85.5

# Update the view state: the window had scrolled so row 15 is at the top,
# and the active selection moved to E10.
$ws.Range("E10").Select() | Out-Null
